$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in newly extracted data for the Bozinovic row (row 7) ---
# New string order matters for shared-string table parity: the "not sure..."
# note is introduced before the "figure 2,3 " note, so set L7 first.
$ws.Range("L7").Value = "not sure how to handle constant mean with changing variance…."
$ws.Range("K7").Value = "figure 2,3 "
$ws.Range("N7").Value = "y"

# --- Fill in newly extracted "extracted" flag for the Delnat row (row 9) ---
$ws.Range("N9").Value = "y"

# --- Correct publication years that were mis-entered ---
# Delnat (row 9): 2016 -> 2019
$ws.Range("F9").Value = 2019
# Haupt (row 10): 2019 -> 2017
$ws.Range("F10").Value = 2017

# --- Remove the Burghardt et al. record entirely (row 8); rows below shift up ---
$ws.Rows(8).Delete()

# --- Refresh the AutoFilter range to the new data extent (A2:T20) ---
$ws.AutoFilterMode = $false
$ws.Range("A2:T20").AutoFilter()

# --- Refresh the recorded sort state to the new data extent (A2:T18) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add2($ws.Range("J2:J18"), 0, 2, 0, 0)
$ws.Sort.SortFields.Add2($ws.Range("O2:O18"), 0, 2, 0, 0)
$ws.Sort.SortFields.Add2($ws.Range("A2:A18"), 0, 1, 0, 0)
$ws.Sort.SetRange($ws.Range("A2:T18"))
$ws.Sort.Apply()

# --- Update the hidden _FilterDatabase defined name to match ---
foreach ($n in $wb.Names) {
    $n.RefersTo = "=Sheet1!`$A`$2:`$T`$20"
}

# --- Move the saved cell selection as recorded in the edited file ---
$ws.Range("N20").Select()
